# Update imputed values produced by RandomForest algorithm re-run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B3"   = 6.171799999999999
    "D3"   = -7.038899999999997
    "A12"  = -21.62310000000001
    "B14"  = 6.361599999999997
    "B26"  = 3.897800000000004
    "D30"  = -7.524900000000001
    "B31"  = 4.871600000000003
    "A32"  = -21.1965
    "B35"  = 9.277800000000003
    "A36"  = -19.7399
    "B37"  = 9.092200000000005
    "A38"  = -19.282
    "D44"  = -7.1649
    "B45"  = 5.191500000000002
    "A46"  = -21.68030000000001
    "A54"  = -21.7474
    "A55"  = -22.5284
    "B57"  = 4.783199999999995
    "D58"  = -8.328999999999997
    "A67"  = -21.44379999999998
    "A69"  = -21.57309999999997
    "A72"  = -21.94410000000001
    "D84"  = -8.853600000000005
    "D89"  = -6.275599999999995
    "A91"  = -21.40110000000001
    "D91"  = -6.225899999999995
    "D92"  = -6.199599999999993
    "A99"  = -20.13999999999998
    "B100" = 5.495299999999999
    "B102" = 8.315
    "D102" = -8.044900000000002
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
